$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0

# Row 2: modality becomes "CBT plus In King", planned assistance bumped to 15000
$ws.Range("F2").Value = "$nbsp" + "CBT plus In King" + "$nbsp "
$ws.Range("H2").Value = 15000

# Row 4: modality switches to "In kind", status flips to ongoing
$ws.Range("F4").Value = "$nbsp" + "In kind" + "$nbsp "
$ws.Range("J4").Value = "ongoing"

# Row 9: modality switches to "In kind", assistance bumped to 15000, status flips to ongoing
$ws.Range("F9").Value = "$nbsp" + "In kind" + "$nbsp "
$ws.Range("H9").Value = 15000
$ws.Range("J9").Value = "ongoing"

# Row 11: status flips to complete
$ws.Range("J11").Value = "complete"

# Update the active selection to match the saved cursor position
$ws.Range("K16").Select()
